$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header-consistent E column values for rows 2-19 -> "Complete"
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 5).Value = "Complete"
}

# Rows 20-25 -> "In progress"
for ($r = 20; $r -le 25; $r++) {
    $ws.Cells.Item($r, 5).Value = "In progress"
}

# Update the selected cell to reflect the new active cell/selection
$ws.Range("G22").Select()
